# Add new power plant types to the Electricity Source subscript on the
# "MCF" (Maximum Capacity Factor) sheet (issues #280 and #99).
#
# New rows 19-24 are appended below the existing table (A1:C18), each
# with a plant-type label in column A and a capacity factor in column B
# that mirrors an existing plant type's value via a same-sheet formula
# reference (matching the source workbook's pattern of "new CCS / new
# technology variants inherit the capacity factor of their base tech").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCF")
$originalActiveSheet = $wb.ActiveSheet
$ws.Activate()

# label, formula-source-cell pairs for the six new rows (19-24)
$newRows = @(
    @{ Row = 19; Label = "hard coal w CCS";                    Ref = "B2"  },
    @{ Row = 20; Label = "natural gas combined cycle w CCS";    Ref = "B4"  },
    @{ Row = 21; Label = "biomass w CCS";                       Ref = "B10" },
    @{ Row = 22; Label = "lignite w CCS";                       Ref = "B14" },
    @{ Row = 23; Label = "small modular reactor";               Ref = "B5"  },
    @{ Row = 24; Label = "hydrogen";                             Ref = "B4"  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Label
    $cell = $ws.Range("B$row")
    $cell.NumberFormat = "0.0000"
    $cell.Formula = "=" + $r.Ref
}

# Match the recorded view state: the next empty row (A25) is left selected
# on the MCF sheet, ready for further data entry.
$ws.Range("A25").Select()

# The workbook was saved while "About" was the focused/active tab; restore
# that so the saved file's tabSelected / activeTab state matches.
$originalActiveSheet.Activate()

Write-Host "Added 6 new plant types (rows 19-24) to MCF sheet"
